# ADC now working; To be done: showing it to the lcd
#
# Typography sheet ("Typography"): row 4 ("Default" typography) gets its
# Wildcard Characters (G) and Wildcard Ranges (H) filled in with "0-9".
#
# Translation sheet ("Translation"): three new translation rows are added
# (rows 4-6), each with a Text ID, the "Default" typography, an alignment,
# the "LTR" direction, and an example translation string for the ADC
# measurement feature.

$wb = $excel.ActiveWorkbook

# --- Typography sheet --------------------------------------------------
$typography = $wb.Worksheets.Item("Typography")
$typography.Range("G4").Value = "0-9"
$typography.Range("H4").Value = "0-9"

# --- Translation sheet --------------------------------------------------
$translation = $wb.Worksheets.Item("Translation")

# Row 4: SingleUseId1 / Default / Left / LTR / "ADC value = <value>"
$translation.Range("B4").Value = "SingleUseId1"
$translation.Range("C4").Value = "Default"
$translation.Range("D4").Value = "Left"
$translation.Range("E4").Value = "LTR"
$translation.Range("F4").Value = "ADC value = <value>"

# Row 5: SingleUseId2 / Default / Left / LTR / "10"
$translation.Range("B5").Value = "SingleUseId2"
$translation.Range("C5").Value = "Default"
$translation.Range("D5").Value = "Left"
$translation.Range("E5").Value = "LTR"
$translation.Range("F5").NumberFormat = "@"
$translation.Range("F5").Value = "10"

# Row 6: SingleUseId3 / Default / Center / LTR / "New ADC"
$translation.Range("B6").Value = "SingleUseId3"
$translation.Range("C6").Value = "Default"
$translation.Range("D6").Value = "Center"
$translation.Range("E6").Value = "LTR"
$translation.Range("F6").Value = "New ADC"
